$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.502.38"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "1.877.05"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7160"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07930"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3104"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08272"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7318"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.876.19"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.287"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "29.488.02"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.915"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007881"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "2.121.52"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.055"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.82%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1626"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.036"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.357"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  +2.57%  "
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05274"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.948"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.200"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7261"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "1.205.40"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.705"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9099"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.08%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "2.016.66"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.793"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.942"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.55%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000121"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.366"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
